# Lay out a small earthquake-analysis time-series table on Sheet1:
# columns A:B hold the "Base" series, columns D:E hold the "Surface" series,
# each with a time (s) / accel (g) header row below its label.
# Order of writes controls shared-string interning order, so the headers
# are written before the series labels to match the expected layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column headers (row 2) for the "Base" table (A:B).
$ws.Range("A2").Value = "time (s)"
$ws.Range("B2").Value = "accel (g)"

# Series labels (row 1).
$ws.Range("D1").Value = "Surface"
$ws.Range("A1").Value = "Base"

# Column headers (row 2) for the "Surface" table (D:E). Columns D/E don't
# carry the A/B column-level number formats, so set them explicitly to
# match the look of the A:B table (0.00 for time, 0.000000 for accel).
$ws.Range("D2").Value = "time (s)"
$ws.Range("D2").NumberFormat = "0.00"
$ws.Range("E2").Value = "accel (g)"
$ws.Range("E2").NumberFormat = "0.000000"

# Leave the selection where the author left it after building the table.
$ws.Range("AC1").Select() | Out-Null
